# Remove the trailing empty paragraph that follows the paragraph
# containing the "_GoBack" bookmark, right before the section break.
# (Commit: "1ª versão - questão 3 inserida" - cleans up a stray blank
# paragraph left at the end of the document.)

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$prevPara = $d.Paragraphs.Item($count - 1)

# Delete from just before the last paragraph mark of the previous
# paragraph through the end of the last (empty) paragraph - this
# removes the final paragraph mark entirely, merging it away so the
# document ends right after the bookmark paragraph.
$deleteRange = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
$deleteRange.Delete()
